# Merge the split "<id>...</id>" runs (produced when the value run used
# RGB color 000000 instead of inheriting the Courier-New/7f6000 style of
# the tag runs) back into a single run, matching the newly downloaded
# tc/tcn/tl content where the whole tag is one run.
#
# Pattern to collapse, found as 3 consecutive runs inside one paragraph:
#   run1: "<id>"      Courier New / 7f6000 / sz18 / szCs18 / rtl0
#   run2: "<value>"   (plain) color 000000 / rtl0
#   run3: "</id>"     Courier New / 7f6000 / sz18 / szCs18 / rtl0
# -> single run: "<id><value></id>" with run1's rPr.

$d = $word.ActiveDocument

$xmlTemplate = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body><w:p><w:r w:rsidDel="00000000" w:rsidR="00000000" w:rsidRPr="00000000">' +
  '<w:rPr>' +
  '<w:rFonts w:ascii="Courier New" w:cs="Courier New" w:eastAsia="Courier New" w:hAnsi="Courier New"/>' +
  '<w:color w:val="7f6000"/>' +
  '<w:sz w:val="18"/>' +
  '<w:szCs w:val="18"/>' +
  '<w:rtl w:val="0"/>' +
  '</w:rPr>' +
  '<w:t xml:space="preserve">{TEXT}</w:t>' +
  '</w:r></w:p></w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

# Walk paragraphs back-to-front so earlier Start/End offsets captured for
# not-yet-processed paragraphs stay valid while we rewrite later ones.
$paraCount = $d.Paragraphs.Count
for ($i = $paraCount; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $pr = $p.Range
    $rawText = $pr.Text
    $text = $rawText.TrimEnd([char]13, [char]7)
    if ($text.Length -lt 10) { continue }
    if ($text.Substring(0,4) -ne "<id>") { continue }
    if ($text.Substring($text.Length - 5) -ne "</id>") { continue }

    $start = $pr.Start
    $openEnd = $start + 4
    $closeStart = $start + $text.Length - 5
    $closeEnd = $start + $text.Length

    $openRun = $d.Range($start, $openEnd)
    $closeRun = $d.Range($closeStart, $closeEnd)
    $midRun = $d.Range($openEnd, $closeStart)

    # Only collapse the specific pattern: "<id>"/"</id>" in Courier New
    # 7f6000 sz18, middle run explicitly colored RGB black (000000) --
    # that's the split-run artifact. Leave runs already uniform (e.g.
    # the "automatic" colored fig_* ids) untouched.
    if ($openRun.Font.Name -ne "Courier New") { continue }
    if ($openRun.Font.Color -ne 24703) { continue }
    if ($closeRun.Font.Name -ne "Courier New") { continue }
    if ($closeRun.Font.Color -ne 24703) { continue }
    if ($midRun.Font.Color -ne 0) { continue }

    $full = $d.Range($start, $closeEnd)
    $fullText = $full.Text

    $escaped = $fullText.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
    $xml = $xmlTemplate.Replace("{TEXT}", $escaped)

    $full.Delete()
    $ins = $d.Range($start, $start)
    $ins.InsertXML($xml)
}

Write-Output "done"
